$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns B:F and I:N for rows 2-25 (pandapower vm_pu results, Case_2_46)
$bfData = @{
    2 = @(1.02,1.055359515953808,1.054647501815008,1.063642249761842,1.072686208689986)
    3 = @(1.02,1.056553590506037,1.055561420612518,1.064880279327619,1.073981476114814)
    4 = @(1.02,1.057325596749082,1.056152156852206,1.065681458810297,1.074819726467995)
    5 = @(1.02,1.057649997371049,1.056400352623156,1.066018298749033,1.075172159286852)
    6 = @(1.02,1.057704456888756,1.056442017037737,1.066074857132099,1.075231336272205)
    7 = @(1.02,1.057329931995656,1.056155473843198,1.065685959584704,1.074824435562534)
    8 = @(1.02,1.055763192953489,1.054956496166032,1.064060629991656,1.07312392571888)
    9 = @(1.02,1.052997391439502,1.052838864211593,1.061197195932286,1.070128270936868)
    10 = @(1.02,1.051150005665538,1.051423752484996,1.059288508666422,1.068131601719875)
    11 = @(1.02,1.050349198328299,1.050810179507036,1.058462058521054,1.067267091933563)
    12 = @(1.02,1.050051607984094,1.050582146321479,1.058155078974964,1.066945980771258)
    13 = @(1.02,1.050115448274438,1.050631065890364,1.058220927136487,1.067014859969159)
    14 = @(1.02,1.050324602186155,1.050791332769377,1.058436683495038,1.067240548673857)
    15 = @(1.02,1.050453450854129,1.050890061888728,1.058569618136005,1.067379603676019)
    16 = @(1.02,1.051203133869877,1.051464455882143,1.059343357710357,1.068188977426372)
    17 = @(1.02,1.051673153229933,1.051824537305219,1.059828708580254,1.068696690508457)
    18 = @(1.02,1.051947223151415,1.052034487675276,1.060111808039182,1.068992837059575)
    19 = @(1.02,1.052040659685284,1.052106061969337,1.060208338147562,1.069093816502961)
    20 = @(1.02,1.051622733356122,1.051785912163623,1.059776634787098,1.068642217091138)
    21 = @(1.02,1.050263015319204,1.050744141627977,1.058373148624284,1.067174088826054)
    22 = @(1.02,1.04940732493435,1.050088415917838,1.057490724625588,1.066251053107844)
    23 = @(1.02,1.049861017720141,1.050436097623047,1.057958514884345,1.066740369487347)
    24 = @(1.02,1.051645516187591,1.051803365446758,1.05980016468497,1.068666831264775)
    25 = @(1.02,1.053713026373843,1.053386908940264,1.061937404804603,1.070902632380209)
}
$inData = @{
    2 = @(1.051309575409007,1.060366788138521,1.057390135443281,1.066360337123209,1.075380064296217,1.061872630603448)
    3 = @(1.051730876426165,1.061211327882237,1.058117081819682,1.067412356293672,1.076490944903978,1.06271836969047)
    4 = @(1.0520018678611,1.061756618285103,1.058586204770663,1.068092609024667,1.07720932780191,1.063264434468254)
    5 = @(1.052115406057348,1.061985576616148,1.058783123656017,1.068378475737859,1.077511234770971,1.063493717946391)
    6 = @(1.052134446974692,1.062024003249474,1.05881616962338,1.068426467533002,1.077561920403571,1.063532199149944)
    7 = @(1.052003386482043,1.061759678742731,1.058588837189292,1.068096429225504,1.077213362292744,1.063267499272084)
    8 = @(1.051452291703544,1.060652450225979,1.057636072003195,1.0667159712475,1.07575558208139,1.062158698363836)
    9 = @(1.050468757410911,1.058692240795913,1.055947472477919,1.064279710228301,1.073183394685336,1.06019570521154)
    10 = @(1.04980465121476,1.057379196793265,1.054815134038804,1.062652912834781,1.071466193777035,1.0588807965358)
    11 = @(1.049515076622696,1.056809132172366,1.054323233826397,1.061947841243745,1.070722026140467,1.058309922357752)
    12 = @(1.049407212222749,1.05659715620466,1.054140279595913,1.061685844960151,1.070445515146605,1.058097645359859)
    13 = @(1.049430363255983,1.056642636146119,1.054179534813664,1.061742048676754,1.070504831989051,1.058143189888051)
    14 = @(1.049506166716011,1.05679161484765,1.054308115690336,1.061926186630245,1.070699171588124,1.058292380156425)
    15 = @(1.049552831476819,1.056883375191636,1.054387306681283,1.062039626602235,1.07081889809516,1.058384270810625)
    16 = @(1.049823826830064,1.057416998160975,1.05484774618897,1.062699692045775,1.071515568693367,1.058918651585788)
    17 = @(1.049993275487081,1.057751320605025,1.055136140802165,1.063113555836984,1.07195240705213,1.059253448806056)
    18 = @(1.050091917952287,1.057946179878751,1.055304203184719,1.063354892307549,1.072207148960362,1.059448584802301)
    19 = @(1.050125519601867,1.058012597262533,1.055361482159372,1.063437171149831,1.072293999521027,1.059515096506387)
    20 = @(1.049975115332529,1.057715466020867,1.05510521466069,1.063069158733196,1.071905544502986,1.059217543304276)
    21 = @(1.049483852874125,1.056747750659746,1.054270258481941,1.061871965409029,1.070641946035808,1.058248453676344)
    22 = @(1.049173220101085,1.056137986063076,1.053743895577509,1.061118655106622,1.069846926378908,1.057637823144006)
    23 = @(1.049338059311963,1.056461359952901,1.054023063098463,1.061518055662928,1.070268433887874,1.057961656261846)
    24 = @(1.049983321724316,1.057731667613211,1.055119189326552,1.06308922006884,1.071926719852553,1.05923376790474)
    25 = @(1.050724503965844,1.059200094125208,1.056385174296063,1.064909995401595,1.073848782007548,1.060704279750815)
}

foreach ($r in $bfData.Keys) {
    $vals = $bfData[$r]
    $arr = New-Object 'object[,]' 1,5
    for ($i = 0; $i -lt 5; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("B$r`:F$r").Value = $arr
}

foreach ($r in $inData.Keys) {
    $vals = $inData[$r]
    $arr = New-Object 'object[,]' 1,6
    for ($i = 0; $i -lt 6; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("I$r`:N$r").Value = $arr
}
